$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-133 down to 87-134
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new Ajo price record
$ws.Cells.Item(86, 1).Value = 7
$ws.Cells.Item(86, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(86, 3).Value = "Ñuble"
$ws.Cells.Item(86, 4).Value = 44488
$ws.Cells.Item(86, 5).Value = 16
$ws.Cells.Item(86, 6).Value = 100112003
$ws.Cells.Item(86, 7).Value = "Ajo"
$ws.Cells.Item(86, 8).Value = "Chino"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 120
$ws.Cells.Item(86, 11).Value = 16000
$ws.Cells.Item(86, 12).Value = 17000
$ws.Cells.Item(86, 13).Value = 16500
$ws.Cells.Item(86, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(86, 15).Value = "China"
$ws.Cells.Item(86, 16).Value = 1650
$ws.Cells.Item(86, 17).Value = 10
$ws.Cells.Item(86, 18).Value = "Hortaliza"
